$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 54

# Copy formatting of the date cell in column A from the row above so the
# new row matches the existing style used throughout the table.
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item($row, 1).Value = 45986
$ws.Cells.Item($row, 2).Value = 2025
$ws.Cells.Item($row, 3).Value = 2.560577522109297
$ws.Cells.Item($row, 4).Value = 2026
$ws.Cells.Item($row, 5).Value = 1.676143333484292
